# Automatische test-sync: 2025-06-23 18:12:50
# Adds a new "Herinnering betaling" row to the Logs sheet, a matching
# "Factuur / Administratie" row to the Dashboard summary sheet, and
# extends the bar chart's category/value series ranges to include it.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# --- Logs!A5:G5 -------------------------------------------------------------
$logs.Range("A5").Value = "Herinnering betaling"
$logs.Range("B5").Value = "mailmind.test@zohomail.eu"
$logs.Range("C5").Value = "Ik zie dat ik nog een openstaande betaling heb. Kunt u dit bevestigen?"
$logs.Range("D5").Value = "Factuur / Administratie"
$logs.Range("E5").Value = "Beste klant,`nDank u voor uw bericht. Om uw openstaande betaling te bevestigen, hebben wij wat meer informatie nodig. Kunt u alstublieft uw factuurnummer en/of klantgegevens doorgeven, zodat wij uw betaling kunnen controleren?`nMet vriendelijke groet,`n[Bedrijfsnaam]"
$logs.Range("F5").Value = "2025-06-23 18:11:50"
$logs.Range("G5").Value = "Ja"

# The new row's wrapped text otherwise leaves the row with an autosized
# (customHeight) row -- re-fit it so it matches the sheet's other rows.
$logs.Rows(5).AutoFit()

# --- Dashboard!A5:B5 ---------------------------------------------------------
$dash.Range("A5").Value = "Factuur / Administratie"
$dash.Range("B5").Value = 1

# --- Extend the existing conditional formatting rules to cover row 5 --------
foreach ($fc in $logs.Range("D2:D4").FormatConditions) {
    $fc.ModifyAppliesToRange($logs.Range("D2:D5"))
}
foreach ($fc in $logs.Range("G2:G4").FormatConditions) {
    $fc.ModifyAppliesToRange($logs.Range("G2:G5"))
}

# --- Extend the chart's source ranges so the new row is plotted -------------
$chart = $dash.ChartObjects(1).Chart
$series = $chart.SeriesCollection(1)
$series.XValues = "='Dashboard'!`$A`$2:`$A`$5"
$series.Values = "='Dashboard'!`$B`$2:`$B`$5"
